$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "CB187"
$ws.Range("B20").Value = 0.73
$ws.Range("C20").Value = 0.93
$ws.Range("D20").Value = 1.02

$ws.Range("A21").Value = "CB194"
$ws.Range("B21").Value = 0.57
$ws.Range("C21").Value = 0.77
$ws.Range("D21").Value = 1.14
